# mcx_d display widget stuff
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: swap a couple of button-mapping numbers (B27/B28, B31/B32) ---
$ws1.Range("B27").Value = 37
$ws1.Range("B28").Value = 39
$ws1.Range("B31").Value = 33
$ws1.Range("B32").Value = 35

# --- Sheet1: new "display widget" reference block in columns Q:T, rows 25-31 ---
# Values are entered in this specific order so the shared-strings table is
# built up in the same sequence as the authored workbook.
$ws1.Range("Q26").Value = "'0x20"
$ws1.Range("Q27").Value = "0x23"
$ws1.Range("R25").Value = "thisaction"
$ws1.Range("R26").Value = "0x2"
$ws1.Range("S25").Value = "thisrow"
$ws1.Range("T25").Value = "thisposition"
$ws1.Range("Q30").Value = "'0x24"
$ws1.Range("Q31").Value = "0x27"

$ws1.Range("R27").Value = "0x2"
$ws1.Range("R30").Value = "0x2"
$ws1.Range("R31").Value = "0x2"

$ws1.Range("S26").Value = 0
$ws1.Range("T26").Value = 0
$ws1.Range("S27").Value = 0
$ws1.Range("T27").Value = 3
$ws1.Range("S30").Value = 1
$ws1.Range("T30").Value = 0
$ws1.Range("S31").Value = 1
$ws1.Range("T31").Value = 3

# Apply left-horizontal alignment style to the new numeric S/T cells (matches style used elsewhere)
$ws1.Range("S26:T27").HorizontalAlignment = -4131
$ws1.Range("S30:T31").HorizontalAlignment = -4131

# --- Sheet view / active tab: Sheet1 becomes the active tab, with a new selection ---
$ws1.Activate() | Out-Null
$ws1.Range("U16").Select() | Out-Null
